# Auto-generated edit script applying the Hades_Profits.xlsx numeric update
# across ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (83 cell(s)) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2818.5278
$ws.Range("I15").Value = 2818.5278
$ws.Range("K15").Value = 8455.5834
$ws.Range("M15").Value = -8286.5834
$ws.Range("H46").Value = 220281.17
$ws.Range("I46").Value = 267062.34
$ws.Range("J46").Value = 173500
$ws.Range("K46").Value = 801187.02
$ws.Range("L46").Value = 520500
$ws.Range("M46").Value = -801068.02
$ws.Range("N46").Value = -520738
$ws.Range("H60").Value = 220281.17
$ws.Range("I60").Value = 267062.34
$ws.Range("J60").Value = 173500
$ws.Range("K60").Value = 801187.02
$ws.Range("L60").Value = 520500
$ws.Range("M60").Value = -800703.02
$ws.Range("N60").Value = -521468
$ws.Range("H62").Value = 2067.742
$ws.Range("I62").Value = 1655.0667
$ws.Range("J62").Value = 2454.625
$ws.Range("K62").Value = 1655.0667
$ws.Range("L62").Value = 2454.625
$ws.Range("M62").Value = -1031.0667
$ws.Range("N62").Value = -3702.625
$ws.Range("H65").Value = 2067.742
$ws.Range("I65").Value = 1655.0667
$ws.Range("J65").Value = 2454.625
$ws.Range("K65").Value = 8275.333500000001
$ws.Range("L65").Value = 12273.125
$ws.Range("M65").Value = -5155.333500000001
$ws.Range("N65").Value = -18513.125
$ws.Range("H98").Value = 947.36365
$ws.Range("I98").Value = 928.5263
$ws.Range("J98").Value = 1066.6666
$ws.Range("K98").Value = 928.5263
$ws.Range("L98").Value = 1066.6666
$ws.Range("M98").Value = 569.4737
$ws.Range("N98").Value = -4062.6666
$ws.Range("H113").Value = 2690.257
$ws.Range("I113").Value = 2581.682
$ws.Range("J113").Value = 2874
$ws.Range("K113").Value = 2581.682
$ws.Range("L113").Value = 2874
$ws.Range("M113").Value = 672.3180000000002
$ws.Range("N113").Value = -9382
$ws.Range("H121").Value = 1408.4375
$ws.Range("J121").Value = 1469
$ws.Range("L121").Value = 4407
$ws.Range("N121").Value = -7901
$ws.Range("H122").Value = 947.36365
$ws.Range("I122").Value = 928.5263
$ws.Range("J122").Value = 1066.6666
$ws.Range("K122").Value = 2785.5789
$ws.Range("L122").Value = 3199.9998
$ws.Range("M122").Value = -335.5789
$ws.Range("N122").Value = -8099.9998
$ws.Range("H129").Value = 851.5294
$ws.Range("J129").Value = 1082.5
$ws.Range("L129").Value = 3247.5
$ws.Range("N129").Value = -13247.5
$ws.Range("H132").Value = 1023650.9
$ws.Range("I132").Value = 3183.7666
$ws.Range("K132").Value = 9551.299800000001
$ws.Range("M132").Value = -7021.299800000001
$ws.Range("H137").Value = 1962351.6
$ws.Range("I137").Value = 2703771
$ws.Range("J137").Value = 2885.7144
$ws.Range("K137").Value = 8111313
$ws.Range("L137").Value = 8657.143199999999
$ws.Range("M137").Value = -8108763
$ws.Range("N137").Value = -13757.1432
$ws.Range("H138").Value = 2009863.5
$ws.Range("I138").Value = 1220.8679
$ws.Range("J138").Value = 5558465.5
$ws.Range("K138").Value = 3662.6037
$ws.Range("L138").Value = 16675396.5
$ws.Range("M138").Value = 1477.3963
$ws.Range("N138").Value = -16685676.5
$ws.Range("H141").Value = 3104.5898
$ws.Range("I141").Value = 1799.0312
$ws.Range("K141").Value = 5397.0936
$ws.Range("M141").Value = -217.0936000000002

# ---- Sheet: ARM (35 cell(s)) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1197.17
$ws.Range("I32").Value = 850.8427
$ws.Range("J32").Value = 3999.2727
$ws.Range("K32").Value = 850.8427
$ws.Range("L32").Value = 3999.2727
$ws.Range("M32").Value = -563.8427
$ws.Range("N32").Value = -4573.2727
$ws.Range("H45").Value = 1150
$ws.Range("I45").Value = 1037.5
$ws.Range("J45").Value = 1600
$ws.Range("K45").Value = 1037.5
$ws.Range("L45").Value = 1600
$ws.Range("M45").Value = -660.5
$ws.Range("N45").Value = -2354
$ws.Range("H61").Value = 18557004
$ws.Range("I61").Value = 21761998
$ws.Range("J61").Value = 128291
$ws.Range("K61").Value = 21761998
$ws.Range("L61").Value = 128291
$ws.Range("M61").Value = -21761786
$ws.Range("N61").Value = -128715
$ws.Range("H122").Value = 2850827.2
$ws.Range("I122").Value = 1869.3334
$ws.Range("J122").Value = 12347354
$ws.Range("K122").Value = 5608.0002
$ws.Range("L122").Value = 37042062
$ws.Range("M122").Value = -3158.0002
$ws.Range("N122").Value = -37046962
$ws.Range("H136").Value = 18557004
$ws.Range("I136").Value = 21761998
$ws.Range("J136").Value = 128291
$ws.Range("K136").Value = 65285994
$ws.Range("L136").Value = 384873
$ws.Range("M136").Value = -65283444
$ws.Range("N136").Value = -389973

# ---- Sheet: CRP (56 cell(s)) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13890680
$ws.Range("I16").Value = 1451.8182
$ws.Range("J16").Value = 25643104
$ws.Range("K16").Value = 1451.8182
$ws.Range("L16").Value = 25643104
$ws.Range("M16").Value = -1164.8182
$ws.Range("N16").Value = -25643678
$ws.Range("H31").Value = 3601.8823
$ws.Range("I31").Value = 1786
$ws.Range("J31").Value = 7960
$ws.Range("K31").Value = 1786
$ws.Range("L31").Value = 7960
$ws.Range("M31").Value = -1491
$ws.Range("N31").Value = -8550
$ws.Range("H34").Value = 3601.8823
$ws.Range("I34").Value = 1786
$ws.Range("J34").Value = 7960
$ws.Range("K34").Value = 1786
$ws.Range("L34").Value = 7960
$ws.Range("M34").Value = -1584
$ws.Range("N34").Value = -8364
$ws.Range("H105").Value = 580
$ws.Range("I105").Value = 580
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 580
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1167
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 13890680
$ws.Range("I113").Value = 1451.8182
$ws.Range("J113").Value = 25643104
$ws.Range("K113").Value = 1451.8182
$ws.Range("L113").Value = 25643104
$ws.Range("M113").Value = 718.1818000000001
$ws.Range("N113").Value = -25647444
$ws.Range("H122").Value = 2247
$ws.Range("I122").Value = 1683.8334
$ws.Range("J122").Value = 2922.8
$ws.Range("K122").Value = 5051.5002
$ws.Range("L122").Value = 8768.400000000001
$ws.Range("M122").Value = -2601.5002
$ws.Range("N122").Value = -13668.4
$ws.Range("H132").Value = 20256.666
$ws.Range("I132").Value = 1511.8462
$ws.Range("J132").Value = 68993.2
$ws.Range("K132").Value = 4535.5386
$ws.Range("L132").Value = 206979.6
$ws.Range("M132").Value = -2005.5386
$ws.Range("N132").Value = -212039.6
$ws.Range("H134").Value = 17367.254
$ws.Range("I134").Value = 1224.3469
$ws.Range("J134").Value = 61311.832
$ws.Range("K134").Value = 3673.0407
$ws.Range("L134").Value = 183935.496
$ws.Range("M134").Value = -1138.0407
$ws.Range("N134").Value = -189005.496

# ---- Sheet: CUL (50 cell(s)) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 188386.83
$ws.Range("I56").Value = 188386.83
$ws.Range("K56").Value = 188386.83
$ws.Range("M56").Value = -187856.83
$ws.Range("H68").Value = 767.9286
$ws.Range("I68").Value = 455.5
$ws.Range("J68").Value = 1184.5
$ws.Range("K68").Value = 1366.5
$ws.Range("L68").Value = 3553.5
$ws.Range("M68").Value = -555.5
$ws.Range("N68").Value = -5175.5
$ws.Range("H71").Value = 767.9286
$ws.Range("I71").Value = 455.5
$ws.Range("J71").Value = 1184.5
$ws.Range("K71").Value = 4099.5
$ws.Range("L71").Value = 10660.5
$ws.Range("M71").Value = -43.5
$ws.Range("N71").Value = -18772.5
$ws.Range("H75").Value = 1334.5
$ws.Range("I75").Value = 1336.1428
$ws.Range("J75").Value = 1332.8572
$ws.Range("K75").Value = 4008.4284
$ws.Range("L75").Value = 3998.5716
$ws.Range("M75").Value = -3010.4284
$ws.Range("N75").Value = -5994.571599999999
$ws.Range("H78").Value = 1334.5
$ws.Range("I78").Value = 1336.1428
$ws.Range("J78").Value = 1332.8572
$ws.Range("K78").Value = 12025.2852
$ws.Range("L78").Value = 11995.7148
$ws.Range("M78").Value = -7033.2852
$ws.Range("N78").Value = -21979.7148
$ws.Range("H131").Value = 1182.1316
$ws.Range("I131").Value = 432.85715
$ws.Range("J131").Value = 1351.3226
$ws.Range("K131").Value = 1298.57145
$ws.Range("L131").Value = 4053.9678
$ws.Range("M131").Value = 3741.42855
$ws.Range("N131").Value = -14133.9678
$ws.Range("H132").Value = 3222.4285
$ws.Range("I132").Value = 2378.8
$ws.Range("J132").Value = 3691.111
$ws.Range("K132").Value = 21409.2
$ws.Range("L132").Value = 33219.999
$ws.Range("M132").Value = -18879.2
$ws.Range("N132").Value = -38279.999
$ws.Range("H141").Value = 8128.7
$ws.Range("I141").Value = 8128.7
$ws.Range("K141").Value = 24386.1
$ws.Range("M141").Value = -19206.1

# ---- Sheet: GSM (18 cell(s)) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1137.68
$ws.Range("I102").Value = 1048.7
$ws.Range("K102").Value = 1048.7
$ws.Range("M102").Value = 573.3
$ws.Range("H122").Value = 1433.0769
$ws.Range("I122").Value = 1259.1428
$ws.Range("J122").Value = 1636
$ws.Range("K122").Value = 3777.4284
$ws.Range("L122").Value = 4908
$ws.Range("M122").Value = -1327.4284
$ws.Range("N122").Value = -9808
$ws.Range("H126").Value = 2245.3125
$ws.Range("I126").Value = 1430.2858
$ws.Range("J126").Value = 2879.2222
$ws.Range("K126").Value = 4290.857400000001
$ws.Range("L126").Value = 8637.6666
$ws.Range("M126").Value = -1820.857400000001
$ws.Range("N126").Value = -13577.6666

# ---- Sheet: LTW (22 cell(s)) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5124.75
$ws.Range("J7").Value = 5714
$ws.Range("L7").Value = 5714
$ws.Range("N7").Value = -5938
$ws.Range("H122").Value = 2945.6765
$ws.Range("I122").Value = 2389.6667
$ws.Range("J122").Value = 3571.1875
$ws.Range("K122").Value = 7169.000100000001
$ws.Range("L122").Value = 10713.5625
$ws.Range("M122").Value = -4719.000100000001
$ws.Range("N122").Value = -15613.5625
$ws.Range("H126").Value = 5124.75
$ws.Range("J126").Value = 5714
$ws.Range("L126").Value = 17142
$ws.Range("N126").Value = -22082
$ws.Range("H132").Value = 30716.143
$ws.Range("I132").Value = 2018.6875
$ws.Range("J132").Value = 54882.42
$ws.Range("K132").Value = 6056.0625
$ws.Range("L132").Value = 164647.26
$ws.Range("M132").Value = -3526.0625
$ws.Range("N132").Value = -169707.26

# ---- Sheet: WVR (25 cell(s)) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1904.3334
$ws.Range("I122").Value = 1334.8
$ws.Range("J122").Value = 2921.3572
$ws.Range("K122").Value = 4004.4
$ws.Range("L122").Value = 8764.071599999999
$ws.Range("M122").Value = -1554.4
$ws.Range("N122").Value = -13664.0716
$ws.Range("H124").Value = 42619.332
$ws.Range("J124").Value = 42619.332
$ws.Range("L124").Value = 42619.332
$ws.Range("N124").Value = -52439.332
$ws.Range("H132").Value = 45591.844
$ws.Range("I132").Value = 30198.53
$ws.Range("J132").Value = 93171.17999999999
$ws.Range("K132").Value = 90595.59
$ws.Range("L132").Value = 279513.54
$ws.Range("M132").Value = -88065.59
$ws.Range("N132").Value = -284573.54
$ws.Range("H136").Value = 42923.613
$ws.Range("I136").Value = 32186.781
$ws.Range("J136").Value = 63134.117
$ws.Range("K136").Value = 96560.34299999999
$ws.Range("L136").Value = 189402.351
$ws.Range("M136").Value = -94010.34299999999
$ws.Range("N136").Value = -194502.351
